$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.063.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.547.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.41"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3826"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3264"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -10.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07297"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.775"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.740"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.556.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001078"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06595"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.359"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.069.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.285"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.01"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.915"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.733.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.073"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.837"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.866"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08190"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.184"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06219"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02301"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.217"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2140"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.232"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5975"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.58"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.720"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5773"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.961"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.33"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.167"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.92%  "
